$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, G across rows 2-12 (F column is unchanged)
$data = @{
    2  = @(0.01253208636536152, 0.04103571897497393, 3.223369029078222, 13.86384647080068, 17.14078330521924)
    3  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    4  = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 2.213936997104367)
    5  = @(0.2881169905109251, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.349763226824225)
    6  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    7  = @(0.2881169905109251, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 16.50004567884073)
    8  = @(0.6545652718822623, 9.983522426115931, 0.7210945179870265, 13.86384647080068, 25.2230286867859)
    9  = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 13.86384647080068, 20.15985084044064)
    10 = @(0.003078177322033415, 0.002658071450198252, 189.6080260415259, 13.86384647080068, 203.4776087610988)
    11 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    12 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 2797.565817734744, 2803.861822104383)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
